$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the header row (row 1) needs a new column inserted at position X (24)
# for the new "flujo_aire" field; existing data rows 2 and 3 are left as-is.
# Shift the header labels from column X (24) through AL (38) one column to
# the right (working right-to-left so values are not overwritten), then set
# the freed-up X1 cell to the new header label.
for ($c = 38; $c -ge 24; $c--) {
    $ws.Cells.Item(1, $c + 1).Value = $ws.Cells.Item(1, $c).Text
}
$ws.Cells.Item(1, 24).Value = "flujo_aire"
